$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule table edits (columns A-C, rows 18-26):
# C18: CONDITION -> ACTION (second rule's condition column repurposed as name header? keep as "ACTION")
$ws.Range("C18").Value = "ACTION"

# B24: was blank -> now holds the new string "gdfvfd"
$ws.Range("B24").Value = "gdfvfd"

# C25 / C26: previously held values, now cleared
$ws.Range("C25").ClearContents()
$ws.Range("C26").ClearContents()

# The second rule pair (columns D:E) is removed entirely across the whole sheet.
$ws.Range("D1:E26").Clear()
